$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) The "_GoBack" bookmark that used to sit right after
#    "git config --list" is gone in the edited document (Word moves
#    this bookmark to track the location of the user's most recent
#    edit). Remove it from its old spot.
# ------------------------------------------------------------------
$removedOld = $false
try {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
    $removedOld = $true
} catch {
    $removedOld = $false
}

# ------------------------------------------------------------------
# 2) Fix the typo "GPA keys" -> "GPG keys" in the sentence
#    "Go to setting, SSH and GPA keys."
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*GPA keys*") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $full = $target.Range.Text
    $pStart = $target.Range.Start
    $idxA = $full.IndexOf("GPA") + 2      # index of the "A" inside "GPA"
    $posA = $pStart + $idxA

    # Replace the single "A" character with "G" ("GPA" -> "GPG").
    $rChar = $d.Range($posA, $posA + 1)
    $rChar.Text = "G"

    # At this point Word's engine merges the edited text back into a
    # single run ("...and GPG" + " keys."). The real document keeps the
    # text "...and GP" and "G" as two separate runs, split exactly where
    # the new "G" was typed. Force that split by dropping a throw-away
    # bookmark on the boundary and removing it again: inserting the
    # bookmark splits the run, and deleting the bookmark afterwards does
    # not re-merge the now-separate runs.
    $splitPos = $posA
    $rSplit = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("_TmpSplit", $rSplit)
    $d.Bookmarks.Item("_TmpSplit").Delete()

    # Re-create the "_GoBack" bookmark exactly where Word's cursor ended
    # up: right after the newly-typed "G", i.e. between "G" and " keys.".
    $bmPos = $posA + 1
    $rBm = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $rBm)
} elseif (-not $removedOld) {
    throw "Could not locate target paragraph containing 'GPA keys' nor the old _GoBack bookmark."
}
